# Building out Python script for harvesting Mac and Linux defaults.
# Added 'Hidden' field to track that separately.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Item" (sheet "Link") table: rename field link_target_id -> target_id
$ws.Range("G4").Value = "target_id"

# "OS_Directory_Item" table: rename field file_id -> item_id
$ws.Range("F5").Value = "item_id"

# "OS_Directory_Item" table: add new field 'hidden'
$ws.Range("F9").Value = "hidden"

# Update the active cell selection to match the author's final cursor position
$ws.Range("F10").Select()
